# Update attendee-count figures ("想去人数") on both the "展览" sheet and
# the aggregated "全部类型" sheet, per the gh-pages regeneration diff.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 13002
$ws1.Range("F9").Value  = 18
$ws1.Range("F10").Value = 12981
$ws1.Range("F13").Value = 8719
$ws1.Range("F14").Value = 7734
$ws1.Range("F15").Value = 205
$ws1.Range("F21").Value = 16
$ws1.Range("F26").Value = 5218

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 13002
$ws4.Range("F10").Value = 18
$ws4.Range("F11").Value = 12981
$ws4.Range("F14").Value = 8719
$ws4.Range("F15").Value = 7734
$ws4.Range("F16").Value = 205
$ws4.Range("F22").Value = 16
$ws4.Range("F29").Value = 5218
